$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. NamedThing: rename header "id" -> "identifier"
# ---------------------------------------------------------------------------
$namedThing = $wb.Worksheets.Item("NamedThing")
$namedThing.Range("A1").Value = "identifier"

# ---------------------------------------------------------------------------
# 2. Dataset: replace the whole header row + drop the vital_status dropdown
#    validation that lived on column D.
# ---------------------------------------------------------------------------
$dataset = $wb.Worksheets.Item("Dataset")
$dataset.Range("D2:D1048576").Validation.Delete()
$dataset.Cells.Clear()

$datasetHeaders = @("authors", "hasPart", "homepage", "keywords", "last-updated", "license", "title", "version", "identifier", "name", "description")
for ($i = 0; $i -lt $datasetHeaders.Length; $i++) {
    $dataset.Cells.Item(1, $i + 1).Value = $datasetHeaders[$i]
}

# ---------------------------------------------------------------------------
# 3. Replace the old "DatasetCollection" sheet with three sheets in the new
#    order: File, Person, DatasetCollection (recreated so it ends up last
#    again, same single "entries" header as before).
# ---------------------------------------------------------------------------
$oldCollection = $wb.Worksheets.Item("DatasetCollection")
[void]$oldCollection.Delete()

$file = $wb.Worksheets.Add($null, $dataset)
$file.Name = "File"
$fileHeaders = @("checksum_md5", "path_posix", "size_in_bytes", "url")
for ($i = 0; $i -lt $fileHeaders.Length; $i++) {
    $file.Cells.Item(1, $i + 1).Value = $fileHeaders[$i]
}

$person = $wb.Worksheets.Add($null, $file)
$person.Name = "Person"
$personHeaders = @("email", "name")
for ($i = 0; $i -lt $personHeaders.Length; $i++) {
    $person.Cells.Item(1, $i + 1).Value = $personHeaders[$i]
}

$collection = $wb.Worksheets.Add($null, $person)
$collection.Name = "DatasetCollection"
$collection.Cells.Item(1, 1).Value = "entries"
